$d = $word.ActiveDocument

# 1. Title text change: NOMBRE -> Conocimientos previos de la arquitectura WIS
$d.Content.Find.Execute("NOMBRE", $false, $false, $false, $false, $false, $true, 1, $false, "Conocimientos previos de la arquitectura WIS", 2)
